$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update numeric values in column C
$ws.Range("C2").Value = 8
$ws.Range("C3").Value = 4
$ws.Range("C5").Value = 11
$ws.Range("C6").Value = 12
$ws.Range("C7").Value = 14
$ws.Range("C8").Value = 11
$ws.Range("C9").Value = 14
$ws.Range("C10").Value = 10
$ws.Range("C11").Value = 4
$ws.Range("C12").Value = 9
$ws.Range("C13").Value = 12
$ws.Range("C14").Value = 14
$ws.Range("C15").Value = 13
$ws.Range("C16").Value = 10
$ws.Range("C17").Value = 13
$ws.Range("C18").Value = 9

# Update text values in column B
$ws.Range("B6").Value = "<mike>"
$ws.Range("B11").Value = "<mike>"
$ws.Range("B14").Value = "<make>"
